$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new row of data (row 5) for LeetCode problem 2880 "Select Data"
$ws.Range("A5").Value = "2880. Select Data"
$ws.Range("B5").Value = "Easy"
$ws.Range("C5").Value = "Data Selecting"
$ws.Range("D5").Value = "Use indexing or .loc. Check pandas.loc: https://pandas.pydata.org/pandas-docs/stable/reference/api/pandas.DataFrame.loc.html "

# Match the style already used on B2:B4 (green fill) for the new Difficulty cell
$ws.Range("B5").Style = $ws.Range("B4").Style

# Add the hyperlink for the Link column (also sets the cell text + hyperlink style)
$ws.Hyperlinks.Add($ws.Range("E5"), "https://leetcode.com/problems/select-data/solutions/4140968/easy-solution-beginner-friendly-pandas-beats-98/?envType=study-plan-v2&envId=introduction-to-pandas&lang=pythondata ", "", "", "https://leetcode.com/problems/select-data/solutions/4140968/easy-solution-beginner-friendly-pandas-beats-98/?envType=study-plan-v2&envId=introduction-to-pandas&lang=pythondata ") | Out-Null

# Grow the existing table (Table2) so it covers the new row
$table = $ws.ListObjects.Item("Table2")
$table.Resize($ws.Range("A1:E5"))

# Update the active selection like the author's session
$ws.Range("E7").Select()
